$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in D/E hold numeric-looking values stored as text (inline strings) in
# the source workbook. Force text formatting before assigning so Excel keeps
# storing them as text instead of silently re-typing them as numbers/percentages.
$updateCells = @("D2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","E21","D22","E22","E23","D24","E24","D25","E25","E26","D27","E27","D40","E40","D41","D42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50")
foreach ($addr in $updateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "270.33"
$ws.Range("D3").Value = "26.72"
$ws.Range("E3").Value = "-1.36%"
$ws.Range("E4").Value = "0.23%"
$ws.Range("D5").Value = "0.06127"
$ws.Range("E5").Value = "-1.30%"
$ws.Range("D6").Value = "6.742"
$ws.Range("E6").Value = "0.24%"
$ws.Range("D7").Value = "0.8538"
$ws.Range("E7").Value = "0.42%"
$ws.Range("D8").Value = "0.8938"
$ws.Range("E8").Value = "-1.90%"
$ws.Range("D9").Value = "0.1426"
$ws.Range("E9").Value = "1.47%"
$ws.Range("D10").Value = "0.05048"
$ws.Range("E10").Value = "7.52%"
$ws.Range("D11").Value = "0.07151"
$ws.Range("E11").Value = "0.80%"
$ws.Range("D12").Value = "0.03164"
$ws.Range("E12").Value = "-0.13%"
$ws.Range("D13").Value = "0.09037"
$ws.Range("E13").Value = "-0.22%"
$ws.Range("D14").Value = "0.001534"
$ws.Range("E14").Value = "0.58%"
$ws.Range("D15").Value = "0.0006081"
$ws.Range("E15").Value = "-1.57%"
$ws.Range("D16").Value = "0.006106"
$ws.Range("E16").Value = "0.44%"
$ws.Range("D17").Value = "3.463"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("D18").Value = "3.179"
$ws.Range("E18").Value = "0.25%"
$ws.Range("E19").Value = "3.98%"
$ws.Range("E21").Value = "-2.21%"
$ws.Range("D22").Value = "3.848"
$ws.Range("E22").Value = "-6.13%"
$ws.Range("E23").Value = "0.41%"
$ws.Range("D24").Value = "0.001178"
$ws.Range("E24").Value = "-3.26%"
$ws.Range("D25").Value = "0.004147"
$ws.Range("E25").Value = "0.45%"
$ws.Range("E26").Value = "0.00%"
$ws.Range("D27").Value = "0.0001681"
$ws.Range("E27").Value = "4.01%"
$ws.Range("D40").Value = "0.03960"
$ws.Range("E40").Value = "1.12%"
$ws.Range("D41").Value = "0.1118"
$ws.Range("D42").Value = "0.004192"
$ws.Range("D43").Value = "0.002038"
$ws.Range("E43").Value = "-6.67%"
$ws.Range("D44").Value = "0.01176"
$ws.Range("E44").Value = "-13.05%"
$ws.Range("D45").Value = "0.00005131"
$ws.Range("E45").Value = "-0.82%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("D47").Value = "0.9040"
$ws.Range("E47").Value = "431.71%"
$ws.Range("D48").Value = "0.02992"
$ws.Range("E48").Value = "-16.67%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "-0.04%"
